$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 10-12 (pushes old rows 10-16 down to 13-19),
# matching formatting from the row above (row 9) and then fixing up.
$ws.Rows("10:12").Insert(-4121, 0)
$ws.Range("A10:P12").ClearFormats()
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Write final label (col B) and A-index for rows 10-19 (data indices 8-17)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.086662413002611
$ws.Range("D10").Value = 0.9491496470760641
$ws.Range("E10").Value = 0.9898537016774713
$ws.Range("F10").Value = 0.9725266603351137
$ws.Range("G10").Value = 1.086662413002611
$ws.Range("H10").Value = 0.9491496470760641
$ws.Range("I10").Value = 1.019375589268239
$ws.Range("J10").Value = 0.9605942201170385
$ws.Range("K10").Value = 1.018040081265707
$ws.Range("L10").Value = 0.9504567340761417
$ws.Range("M10").Value = 1.086662413002611
$ws.Range("N10").Value = 0.9695016743767677
$ws.Range("O10").Value = 0.999548105522815
$ws.Range("P10").Value = 0.9933323808522984

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.5875720723398838
$ws.Range("D11").Value = 1.534142505545308
$ws.Range("E11").Value = 0.905430746801493
$ws.Range("F11").Value = 1.097881806962828
$ws.Range("G11").Value = 0.5875720723398838
$ws.Range("H11").Value = 1.534142505545308
$ws.Range("I11").Value = 0.8073493666100912
$ws.Range("J11").Value = 1.117987974520836
$ws.Range("K11").Value = 0.8391002428608455
$ws.Range("L11").Value = 1.342545656118851
$ws.Range("M11").Value = 0.5875720723398838
$ws.Range("N11").Value = 1.219786626173401
$ws.Range("O11").Value = 1.031256782912378
$ws.Range("P11").Value = 1.029001296470017

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.586606018181125
$ws.Range("D12").Value = 1.536670656902683
$ws.Range("E12").Value = 0.9042561087401912
$ws.Range("F12").Value = 1.098484798449754
$ws.Range("G12").Value = 0.586606018181125
$ws.Range("H12").Value = 1.536670656902683
$ws.Range("I12").Value = 0.8062069664315179
$ws.Range("J12").Value = 1.118218474851072
$ws.Range("K12").Value = 0.8386816357067002
$ws.Range("L12").Value = 1.344255264380744
$ws.Range("M12").Value = 0.586606018181125
$ws.Range("N12").Value = 1.220463382821437
$ws.Range("O12").Value = 1.031504395568438
$ws.Range("P12").Value = 1.029172490455474

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.5873699306972402
$ws.Range("D13").Value = 1.534528374753347
$ws.Range("E13").Value = 0.9052205473826705
$ws.Range("F13").Value = 1.098030750560502
$ws.Range("G13").Value = 0.5873699306972402
$ws.Range("H13").Value = 1.534528374753347
$ws.Range("I13").Value = 0.8070939819007268
$ws.Range("J13").Value = 1.11803308165379
$ws.Range("K13").Value = 0.839023427125616
$ws.Range("L13").Value = 1.342900700843288
$ws.Range("M13").Value = 0.5873699306972402
$ws.Range("N13").Value = 1.219874461068009
$ws.Range("O13").Value = 1.03128740084844
$ws.Range("P13").Value = 1.029025099364647

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.7451279999999989
$ws.Range("D14").Value = 1.105512000000001
$ws.Range("E14").Value = 1.107011999999999
$ws.Range("F14").Value = 0.9982560000000004
$ws.Range("G14").Value = 0.7451279999999989
$ws.Range("H14").Value = 1.105512000000001
$ws.Range("I14").Value = 0.9840999999999996
$ws.Range("J14").Value = 1.077896
$ws.Range("K14").Value = 0.9013039999999988
$ws.Range("L14").Value = 1.086659999999999
$ws.Range("M14").Value = 0.7451279999999989
$ws.Range("N14").Value = 1.106262
$ws.Range("O14").Value = 0.9889769999999998
$ws.Range("P14").Value = 1.000733499999999

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.96
$ws.Range("D15").Value = 0.53
$ws.Range("E15").Value = 1.36
$ws.Range("F15").Value = 0.8711000000000011
$ws.Range("G15").Value = 0.96
$ws.Range("H15").Value = 0.53
$ws.Range("I15").Value = 1.220550000000002
$ws.Range("J15").Value = 1.03
$ws.Range("K15").Value = 0.99
$ws.Range("L15").Value = 0.7298624999999997
$ws.Range("M15").Value = 0.96
$ws.Range("N15").Value = 0.9450000000000001
$ws.Range("O15").Value = 0.9302750000000003
$ws.Range("P15").Value = 0.9614390625000004

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9778469556223971
$ws.Range("D16").Value = 0.7241934102528013
$ws.Range("E16").Value = 1.204740374527997
$ws.Range("F16").Value = 0.9240014787584009
$ws.Range("G16").Value = 0.9778469556223971
$ws.Range("H16").Value = 0.7241934102528013
$ws.Range("I16").Value = 1.123437328588797
$ws.Range("J16").Value = 1.013053178777598
$ws.Range("K16").Value = 0.9929532376064024
$ws.Range("L16").Value = 0.8407857504256011
$ws.Range("M16").Value = 0.9778438189055974
$ws.Range("N16").Value = 0.9644668923903992
$ws.Range("O16").Value = 0.957695554790399
$ws.Range("P16").Value = 0.9751264643199993

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9944007598602194
$ws.Range("D17").Value = 0.9943698965328298
$ws.Range("E17").Value = 0.9946035920634615
$ws.Range("F17").Value = 0.9938119576426394
$ws.Range("G17").Value = 0.9944007598602194
$ws.Range("H17").Value = 0.9943698965328298
$ws.Range("I17").Value = 0.9938804113612801
$ws.Range("J17").Value = 0.9951694324111999
$ws.Range("K17").Value = 0.9940530526952971
$ws.Range("L17").Value = 0.9930084784487491
$ws.Range("M17").Value = 0.9943838929195346
$ws.Range("N17").Value = 0.9944867442981457
$ws.Range("O17").Value = 0.9942965515247876
$ws.Range("P17").Value = 0.9941621976269595

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9916620388927332
$ws.Range("D18").Value = 1.0121896884806
$ws.Range("E18").Value = 0.9871724681607575
$ws.Range("F18").Value = 0.9955899502848872
$ws.Range("G18").Value = 0.9916620388927332
$ws.Range("H18").Value = 1.0121896884806
$ws.Range("I18").Value = 0.9897174245772615
$ws.Range("J18").Value = 0.993945994030054
$ws.Range("K18").Value = 0.9930499070793373
$ws.Range("L18").Value = 1.003020687710952
$ws.Range("M18").Value = 0.9916620388927332
$ws.Range("N18").Value = 0.9996810783206789
$ws.Range("O18").Value = 0.9966535364547444
$ws.Range("P18").Value = 0.9957935199020728

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9751891190894817
$ws.Range("D19").Value = 1.035361777999287
$ws.Range("E19").Value = 0.9828187211147899
$ws.Range("F19").Value = 1.002396193712235
$ws.Range("G19").Value = 0.9751891190894817
$ws.Range("H19").Value = 1.035361777999287
$ws.Range("I19").Value = 0.9777807969986261
$ws.Range("J19").Value = 0.9984958748600403
$ws.Range("K19").Value = 0.9859507265059233
$ws.Range("L19").Value = 1.020063677084398
$ws.Range("M19").Value = 0.9751552076613681
$ws.Range("N19").Value = 1.009090249557038
$ws.Range("O19").Value = 0.9989414529789484
$ws.Range("P19").Value = 0.9972571109205977
